# Applies the cryptocurrency price/volume refresh described in the
# "Updated cryptos list on Wed Sep  4 17:51:48 UTC 2024 with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '58.027.55'
$ws.Range('E2').Value = '  +0.25%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '2.415.01'
$ws.Range('E3').Value = '  -1.37%  '

# Row 4: TetherUSD
$ws.Range('D4').Value = '''0.999'
$ws.Range('D4').Style = 'Normal'  # drop quote-prefix style, keep cell as plain text
$ws.Range('E4').Value = '  -0.13%  '

# Row 5: BNB
$ws.Range('D5').Value = '''509.07'
$ws.Range('D5').Style = 'Normal'  # drop quote-prefix style, keep cell as plain text

# Row 6: Solana
$ws.Range('D6').Value = '''132.96'
$ws.Range('D6').Style = 'Normal'  # drop quote-prefix style, keep cell as plain text
$ws.Range('E6').Value = '  +1.57%  '

# Row 7: USDC
$ws.Range('D7').Value = '''0.995'
$ws.Range('D7').Style = 'Normal'  # drop quote-prefix style, keep cell as plain text
$ws.Range('E7').Value = '  -0.52%  '

# Row 8: XRP
$ws.Range('E8').Value = '  -0.62%  '

# Row 9: LidoStakedEther
$ws.Range('D9').Value = '2.452.68'
$ws.Range('E9').Value = '  +0.05%  '

# Row 10: Dogecoin
$ws.Range('D10').Value = '''0.0988'
$ws.Range('D10').Style = 'Normal'  # drop quote-prefix style, keep cell as plain text
$ws.Range('E10').Value = '  +0.39%  '

# Row 11: TRON
$ws.Range('E11').Value = '  -1.37%  '

# Row 12: Cardano
$ws.Range('E12').Value = '  -0.36%  '

# Row 13: Toncoin
$ws.Range('D13').Value = '''4.65'
$ws.Range('D13').Style = 'Normal'  # drop quote-prefix style, keep cell as plain text
$ws.Range('E13').Value = '  -5.67%  '

# Row 14: WrappedliquidstakedEther2.0
$ws.Range('D14').Value = '2.849.09'
$ws.Range('E14').Value = '  -1.28%  '

# Row 15: WrappedBTC
$ws.Range('D15').Value = '57.453.39'
$ws.Range('E15').Value = '  -0.60%  '

# Row 16: Avalanche
$ws.Range('D16').Value = '''21.99'
$ws.Range('D16').Style = 'Normal'  # drop quote-prefix style, keep cell as plain text
$ws.Range('E16').Value = '  +1.35%  '

# Row 17: ShibaInu
$ws.Range('E17').Value = '  +0.75%  '

# Row 18: WrappedEther
$ws.Range('D18').Value = '2.434.79'
$ws.Range('E18').Value = '  -0.40%  '

# Row 19: Chainlink
$ws.Range('D19').Value = '''10.35'
$ws.Range('D19').Style = 'Normal'  # drop quote-prefix style, keep cell as plain text
$ws.Range('E19').Value = '  +0.53%  '

# Row 20: Polkadot
$ws.Range('E20').Value = '  +0.00%  '

# Row 21: BitcoinCash
$ws.Range('D21').Value = '''315.19'
$ws.Range('D21').Style = 'Normal'  # drop quote-prefix style, keep cell as plain text
$ws.Range('E21').Value = '  -0.49%  '

# Row 22: Uniswap
$ws.Range('D22').Value = '''6.38'
$ws.Range('D22').Style = 'Normal'  # drop quote-prefix style, keep cell as plain text
$ws.Range('E22').Value = '  +4.52%  '

# Row 23: Dai
$ws.Range('E23').Value = '  -0.40%  '

# Row 24: LEO
$ws.Range('D24').Value = '''5.69'
$ws.Range('D24').Style = 'Normal'  # drop quote-prefix style, keep cell as plain text
$ws.Range('E24').Value = '  -2.32%  '

# Row 25: Litecoin
$ws.Range('D25').Value = '''65.46'
$ws.Range('D25').Style = 'Normal'  # drop quote-prefix style, keep cell as plain text
$ws.Range('E25').Value = '  +1.35%  '

# Row 26: WrappedeETH
$ws.Range('B26').Value = 'Binance-PegBSC-USD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D26').Value = '''0.994'
$ws.Range('D26').Style = 'Normal'  # drop quote-prefix style, keep cell as plain text
$ws.Range('E26').Value = '  -0.55%  '

# Row 27: Binance-PegBSC-USD
$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').Value = '2.548.99'

# Row 28: Polygon
$ws.Range('E28').Value = '  -4.46%  '

# Row 29: Kaspa
$ws.Range('E29').Value = '  -1.88%  '

# Row 30: InternetComputer(DFINITY)
$ws.Range('D30').Value = '''7.66'
$ws.Range('D30').Style = 'Normal'  # drop quote-prefix style, keep cell as plain text
$ws.Range('E30').Value = '  +5.81%  '

# Row 31: Monero
$ws.Range('D31').Value = '''174.21'
$ws.Range('D31').Style = 'Normal'  # drop quote-prefix style, keep cell as plain text
$ws.Range('E31').Value = '  -0.04%  '

# Row 32: PEPE
$ws.Range('D32').Value = '0.0₃0741'
$ws.Range('E32').Value = '  +0.57%  '

# Row 33: PancakeSwap
$ws.Range('E33').Value = '  -0.13%  '

# Row 34: Aptos
$ws.Range('D34').Value = '''6.24'
$ws.Range('D34').Style = 'Normal'  # drop quote-prefix style, keep cell as plain text
$ws.Range('E34').Value = '  +1.85%  '

# Row 35: Fetch.AI
$ws.Range('E35').Value = '  +0.30%  '

# Row 37: FirstDigitalUSD
$ws.Range('D37').Value = '''0.994'
$ws.Range('D37').Style = 'Normal'  # drop quote-prefix style, keep cell as plain text
$ws.Range('E37').Value = '  -0.32%  '

# Row 38: EthereumClassic
$ws.Range('D38').Value = '''18.05'
$ws.Range('D38').Style = 'Normal'  # drop quote-prefix style, keep cell as plain text
$ws.Range('E38').Value = '  +1.29%  '

# Row 39: ImmutableX
$ws.Range('E39').Value = '  +5.20%  '

# Row 40: NEARProtocol
$ws.Range('D40').Value = '''3.90'
$ws.Range('D40').Style = 'Normal'  # drop quote-prefix style, keep cell as plain text
$ws.Range('E40').Value = '  +3.49%  '

# Row 41: SuiNetwork
$ws.Range('E41').Value = '  +3.60%  '

# Row 42: OKB
$ws.Range('D42').Value = '''36.54'
$ws.Range('D42').Style = 'Normal'  # drop quote-prefix style, keep cell as plain text
$ws.Range('E42').Value = '  +0.82%  '

# Row 43: Stacks
$ws.Range('D43').Value = '''1.48'
$ws.Range('D43').Style = 'Normal'  # drop quote-prefix style, keep cell as plain text
$ws.Range('E43').Value = '  +1.79%  '

# Row 44: Aave
$ws.Range('D44').Value = '''134.60'
$ws.Range('D44').Style = 'Normal'  # drop quote-prefix style, keep cell as plain text
$ws.Range('E44').Value = '  +9.48%  '

# Row 45: Filecoin
$ws.Range('D45').Value = '''3.43'
$ws.Range('D45').Style = 'Normal'  # drop quote-prefix style, keep cell as plain text
$ws.Range('E45').Value = '  +0.27%  '

# Row 46: RenderToken
$ws.Range('E46').Value = '  +4.68%  '

# Row 47: Bittensor
$ws.Range('D47').Value = '''260.64'
$ws.Range('D47').Style = 'Normal'  # drop quote-prefix style, keep cell as plain text
$ws.Range('E47').Value = '  -1.13%  '

# Row 48: Mantle
$ws.Range('E48').Value = '  -2.28%  '

# Row 49: Stellar
$ws.Range('E49').Value = '  -0.41%  '

# Row 50: Hedera
$ws.Range('D50').Value = '''0.0497'
$ws.Range('D50').Style = 'Normal'  # drop quote-prefix style, keep cell as plain text
$ws.Range('E50').Value = '  +0.33%  '

# Row 51: VeChain
$ws.Range('E51').Value = '  +1.68%  '
